$ws = $excel.ActiveWorkbook.ActiveSheet

$cells = [ordered]@{
    'D2' = '26.588.13'
    'E2' = '  -0.16%  '
    'D3' = '1.730.42'
    'E3' = '  -1.00%  '
    'D4' = '0.9994'
    'E4' = '  +0.00%  '
    'D5' = '246.05'
    'E5' = '  -0.74%  '
    'D6' = '0.9997'
    'E6' = '  -0.05%  '
    'E7' = '  +0.09%  '
    'D8' = '0.2669'
    'E8' = '  -1.19%  '
    'D9' = '0.06180'
    'E9' = '  -1.34%  '
    'D10' = '1.732.24'
    'E10' = '  -0.66%  '
    'D11' = '0.07114'
    'E11' = '  +0.01%  '
    'D12' = '15.63'
    'E12' = '  -1.16%  '
    'D13' = '0.6126'
    'E13' = '  -0.82%  '
    'D14' = '4.548'
    'E14' = '  +0.77%  '
    'D15' = '77.32'
    'E15' = '  +0.03%  '
    'D16' = '0.9995'
    'E16' = '  -0.07%  '
    'D17' = '26.589.28'
    'E17' = '  -0.14%  '
    'D18' = '0.9998'
    'E18' = '  -0.04%  '
    'D19' = '0.000006969'
    'E19' = '  +0.83%  '
    'D20' = '11.55'
    'E20' = '  -1.48%  '
    'D21' = '1.954.27'
    'E21' = '  -0.83%  '
    'D22' = '4.525'
    'E22' = '  -2.77%  '
    'E23' = '  -0.54%  '
    'D24' = '5.244'
    'E24' = '  -2.15%  '
    'D25' = '137.25'
    'E25' = '  +0.77%  '
    'D26' = '15.37'
    'E26' = '  -0.68%  '
    'D27' = '1.782'
    'E27' = '  -2.34%  '
    'E28' = '  +0.03%  '
    'D29' = '108.42'
    'E29' = '  +0.55%  '
    'D30' = '3.977'
    'E30' = '  -1.28%  '
    'D31' = '0.08019'
    'E31' = '  +1.50%  '
    'D32' = '3.689'
    'E32' = '  -2.37%  '
    'D33' = '0.04539'
    'E33' = '  -1.00%  '
    'B34' = 'HuobiToken'
    'C34' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D34' = '2.616'
    'E34' = '  +0.03%  '
    'B35' = 'ARBITRUM'
    'C35' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D35' = '1.004'
    'E35' = '  +0.46%  '
    'B36' = 'ImmutableX'
    'C36' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D36' = '0.6342'
    'E36' = '  -0.22%  '
    'B37' = 'RenderToken'
    'C37' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D37' = '2.055'
    'E37' = '  +3.54%  '
    'B38' = 'TrustWalletToken'
    'C38' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D38' = '0.8980'
    'E38' = '  -5.54%  '
    'B39' = 'MXToken'
    'C39' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D39' = '2.390'
    'E39' = '  -3.07%  '
    'B40' = 'PaxDollar'
    'C40' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D40' = '1.004'
    'E40' = '  -0.03%  '
    'B41' = 'Quant'
    'C41' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D41' = '102.92'
    'E41' = '  -10.09%  '
    'B42' = 'VeChain'
    'C42' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D42' = '0.01501'
    'E42' = '  -0.79%  '
    'B43' = 'FraxShare'
    'C43' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D43' = '5.458'
    'E43' = '  -4.16%  '
    'B44' = 'Aptos'
    'C44' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D44' = '7.219'
    'E44' = '  +6.93%  '
    'B45' = 'TheSandbox'
    'C45' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D45' = '0.3904'
    'E45' = '  -0.44%  '
    'B46' = 'Algorand'
    'C46' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D46' = '0.1186'
    'E46' = '  -1.49%  '
    'B47' = 'Cronos'
    'C47' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D47' = '0.05388'
    'E47' = '  +1.12%  '
    'B48' = 'EnergySwap'
    'C48' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D48' = '7.891'
    'E48' = '  -1.20%  '
    'B49' = 'Elrond'
    'C49' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
    'D49' = '30.70'
    'B50' = 'NEARProtocol'
    'C50' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D50' = '1.256'
    'E50' = '  -0.03%  '
    'B51' = 'Decentraland'
    'C51' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'D51' = '0.3415'
    'E51' = '  -1.27%  '
}

foreach ($ref in $cells.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$ref]
    $rng.NumberFormat = "General"
}
